$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$styleObj = $ws.Range("B16").Style
$ws.Range("Z1").Style = $styleObj
Write-Host "Z1 style:" $ws.Range("Z1").Style.Name
Write-Host "done"
